$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 -> DATA ---
$runmanager = $wb.Worksheets.Item("RUNMANAGER")
$dataSheet  = $wb.Worksheets.Item("Sheet2")
$dataSheet.Name = "DATA"

# --- RUNMANAGER: add a new row for the Amazon hamburger-menu test ---
$runmanager.Range("A4").Value = "amazonHamburgerMenuTest"
$runmanager.Range("B4").Value = "To check whether Amazon website is working or not"
$runmanager.Range("C4").Value = "yes"
$runmanager.Range("D4").Value = "1"
$runmanager.Range("E4").Value = "1"

$runmanager.Columns.Item(1).ColumnWidth = 27.140625

$runmanager.Range("A4").Select()

# --- DATA sheet: populate the data-provider table ---
$header = @("testcasename","execute","browser","username","password","name","menuItem")
for ($c = 0; $c -lt $header.Length; $c++) {
    $dataSheet.Cells.Item(1, $c + 1).Value = $header[$c]
}

$rows = @(
    @("loginLogoutTest","yes","chrome","Admin","admin123","Niyaz",""),
    @("loginLogoutTest","yes","edge","Admin","admin123","Subscribe",""),
    @("newTest","yes","chrome","Admin","admin123","",""),
    @("newTest","yes","edge","Admin","admin123","",""),
    @("loginLogoutTest","yes","edge","admin123","admin123","SeleniumAutomation",""),
    @("amazonHamburgerMenuTest","yes","edge","","","","Power Banks")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $dataSheet.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

$dataSheet.Columns.Item(1).ColumnWidth = 27.140625
$dataSheet.Columns.Item(4).ColumnWidth = 20.85546875
$dataSheet.Columns.Item(5).ColumnWidth = 15
$dataSheet.Columns.Item(6).ColumnWidth = 20.7109375
$dataSheet.Columns.Item(7).ColumnWidth = 19.140625

$dataSheet.Range("C7").Select()
$dataSheet.Activate()
